$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("F2").Value = "2021-10-05 13:39:26.158244"
$ws.Range("F3").Value = "2021-10-05 13:39:26.158255"
$ws.Range("F4").Value = "2021-10-05 13:39:26.158259"
$ws.Range("F5").Value = "2021-10-05 13:39:26.158261"
$ws.Range("F6").Value = "2021-10-05 13:39:26.158264"
$ws.Range("F7").Value = "2021-10-05 13:39:26.158267"
$ws.Range("F8").Value = "2021-10-05 13:39:26.158270"
$ws.Range("F9").Value = "2021-10-05 13:39:26.158272"
$ws.Range("F10").Value = "2021-10-05 13:39:26.158275"
$ws.Range("F11").Value = "2021-10-05 13:39:26.158278"
$ws.Range("F12").Value = "2021-10-05 13:39:26.158280"
$ws.Range("F13").Value = "2021-10-05 13:39:26.158283"
$ws.Range("F14").Value = "2021-10-05 13:39:26.158285"
$ws.Range("F15").Value = "2021-10-05 13:39:26.158288"
$ws.Range("F16").Value = "2021-10-05 13:39:26.158291"
$ws.Range("F17").Value = "2021-10-05 13:39:26.158293"
$ws.Range("F18").Value = "2021-10-05 13:39:26.158296"
$ws.Range("F19").Value = "2021-10-05 13:39:26.158298"
$ws.Range("F20").Value = "2021-10-05 13:39:26.158301"
$ws.Range("F21").Value = "2021-10-05 13:39:26.158304"
$ws.Range("F22").Value = "2021-10-05 13:39:26.158306"
$ws.Range("F23").Value = "2021-10-05 13:39:26.158309"
$ws.Range("F24").Value = "2021-10-05 13:39:26.158311"
$ws.Range("F25").Value = "2021-10-05 13:39:26.158314"
$ws.Range("F26").Value = "2021-10-05 13:39:26.158316"
$ws.Range("F27").Value = "2021-10-05 13:39:26.158319"
$ws.Range("F28").Value = "2021-10-05 13:39:26.158322"
$ws.Range("F29").Value = "2021-10-05 13:39:26.158324"
$ws.Range("F30").Value = "2021-10-05 13:39:26.158327"
$ws.Range("F31").Value = "2021-10-05 13:39:26.158329"
$ws.Range("F32").Value = "2021-10-05 13:39:26.158332"
$ws.Range("F33").Value = "2021-10-05 13:39:26.158334"
$ws.Range("F34").Value = "2021-10-05 13:39:26.158337"
$ws.Range("F35").Value = "2021-10-05 13:39:26.158340"
$ws.Range("F36").Value = "2021-10-05 13:39:26.158342"
$ws.Range("F37").Value = "2021-10-05 13:39:26.158345"
$ws.Range("F38").Value = "2021-10-05 13:39:26.158347"
$ws.Range("F39").Value = "2021-10-05 13:39:26.158350"
$ws.Range("F40").Value = "2021-10-05 13:39:26.158352"
$ws.Range("F41").Value = "2021-10-05 13:39:26.158355"
$ws.Range("F42").Value = "2021-10-05 13:39:26.158358"
$ws.Range("F43").Value = "2021-10-05 13:39:26.158360"
$ws.Range("F44").Value = "2021-10-05 13:39:26.158363"
$ws.Range("F45").Value = "2021-10-05 13:39:26.158365"
$ws.Range("F46").Value = "2021-10-05 13:39:26.158368"
$ws.Range("F47").Value = "2021-10-05 13:39:26.158371"
$ws.Range("F48").Value = "2021-10-05 13:39:26.158373"
$ws.Range("F49").Value = "2021-10-05 13:39:26.158376"
$ws.Range("F50").Value = "2021-10-05 13:39:26.158378"
$ws.Range("F51").Value = "2021-10-05 13:39:26.158381"
$ws.Range("F52").Value = "2021-10-05 13:39:26.158383"
$ws.Range("F53").Value = "2021-10-05 13:39:26.158386"
$ws.Range("F54").Value = "2021-10-05 13:39:26.158389"
$ws.Range("F55").Value = "2021-10-05 13:39:26.158392"
